# Update "想去人数" (column F) values on sheets "展览" and "全部类型"
# as published in the newer site snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# row -> new value for column F
$updates = @{
    2  = 1079
    3  = 784
    8  = 1921
    9  = 7112
    10 = 893
    11 = 394
    12 = 328
    13 = 117
    14 = 389
    16 = 7049
    17 = 287
    18 = 1318
    19 = 143
    22 = 124
    23 = 290
    24 = 124
    26 = 14
    28 = 17
    30 = 601
    31 = 52
    32 = 87
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
